# Auto-generated script applying the diff to Zalera_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 449
$ws.Range("I12").Value = 449
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 449
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -279
$ws.Range("N12").ClearContents()

$ws.Range("H18").Value = 725.2
$ws.Range("I18").Value = 725.2
$ws.Range("K18").Value = 725.2
$ws.Range("M18").Value = -441.2

$ws.Range("H38").Value = 139.3
$ws.Range("I38").Value = 139.3
$ws.Range("K38").Value = 417.9
$ws.Range("M38").Value = -45.90000000000003

$ws.Range("H43").Value = 15661
$ws.Range("J43").Value = 11966.083
$ws.Range("L43").Value = 11966.083
$ws.Range("N43").Value = -12104.083

$ws.Range("H47").Value = 51424.5
$ws.Range("I47").Value = 14849.5
$ws.Range("K47").Value = 14849.5
$ws.Range("M47").Value = -13877.5

$ws.Range("H62").Value = 33340334
$ws.Range("I62").Value = 33340334
$ws.Range("K62").Value = 33340334
$ws.Range("M62").Value = -33339710

$ws.Range("H65").Value = 33340334
$ws.Range("I65").Value = 33340334
$ws.Range("K65").Value = 166701670
$ws.Range("M65").Value = -166698550

$ws.Range("H98").Value = 1140.3077
$ws.Range("J98").Value = 1594.25
$ws.Range("L98").Value = 1594.25
$ws.Range("N98").Value = -4590.25

$ws.Range("H107").Value = 19232056
$ws.Range("I107").Value = 22728216
$ws.Range("J107").Value = 3171.5
$ws.Range("K107").Value = 22728216
$ws.Range("L107").Value = 3171.5
$ws.Range("M107").Value = -22726296
$ws.Range("N107").Value = -7011.5

$ws.Range("H113").Value = 170168
$ws.Range("I113").Value = 2667.6667
$ws.Range("J113").Value = 337668.34
$ws.Range("K113").Value = 2667.6667
$ws.Range("L113").Value = 337668.34
$ws.Range("M113").Value = 586.3332999999998
$ws.Range("N113").Value = -344176.34

$ws.Range("H122").Value = 1140.3077
$ws.Range("J122").Value = 1594.25
$ws.Range("L122").Value = 4782.75
$ws.Range("N122").Value = -9682.75

$ws.Range("H125").Value = 1499
$ws.Range("I125").Value = 1499
$ws.Range("K125").Value = 13491
$ws.Range("M125").Value = -11031

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2246573
$ws.Range("I2").Value = 2246573
$ws.Range("K2").Value = 2246573
$ws.Range("M2").Value = -2246460

$ws.Range("H32").Value = 38555.91
$ws.Range("I32").Value = 42410.793
$ws.Range("K32").Value = 42410.793
$ws.Range("M32").Value = -42123.793

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H116").Value = 2246573
$ws.Range("I116").Value = 2246573
$ws.Range("K116").Value = 2246573
$ws.Range("M116").Value = -2244279

$ws.Range("H122").Value = 3354.5386
$ws.Range("I122").Value = 2283.1667
$ws.Range("J122").Value = 4272.857
$ws.Range("K122").Value = 6849.500100000001
$ws.Range("L122").Value = 12818.571
$ws.Range("M122").Value = -4399.500100000001
$ws.Range("N122").Value = -17718.571

$ws.Range("H132").Value = 3417.0352
$ws.Range("I132").Value = 2632.1914
$ws.Range("K132").Value = 7896.574200000001
$ws.Range("M132").Value = -5366.574200000001

$ws.Range("H139").Value = 99699.5
$ws.Range("J139").Value = 98798
$ws.Range("L139").Value = 98798
$ws.Range("N139").Value = -109078

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2246573
$ws.Range("I3").Value = 2246573
$ws.Range("K3").Value = 2246573
$ws.Range("M3").Value = -2246459

$ws.Range("H22").Value = 3425.5
$ws.Range("I22").Value = 3425.5
$ws.Range("K22").Value = 3425.5
$ws.Range("M22").Value = -3252.5

$ws.Range("H94").Value = 2093.9333
$ws.Range("I94").Value = 1999.875
$ws.Range("K94").Value = 1999.875
$ws.Range("M94").Value = -1548.875

$ws.Range("H107").Value = 1381.4117
$ws.Range("I107").Value = 1305.0769
$ws.Range("K107").Value = 1305.0769
$ws.Range("M107").Value = 614.9231

$ws.Range("H124").Value = 79499.5
$ws.Range("J124").Value = 79499.5
$ws.Range("L124").Value = 79499.5
$ws.Range("N124").Value = -89319.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H132").Value = 4209.2666
$ws.Range("I132").Value = 2356.8462
$ws.Range("K132").Value = 7070.5386
$ws.Range("M132").Value = -4540.5386

$ws.Range("H134").Value = 3480.923
$ws.Range("I134").Value = 2560.5833
$ws.Range("K134").Value = 7681.749899999999
$ws.Range("M134").Value = -5146.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 201
$ws.Range("I7").Value = 170.125
$ws.Range("J7").Value = 283.33334
$ws.Range("K7").Value = 510.375
$ws.Range("L7").Value = 850.0000200000001
$ws.Range("M7").Value = -398.375
$ws.Range("N7").Value = -1074.00002

$ws.Range("H12").Value = 417.5
$ws.Range("J12").Value = 417.5
$ws.Range("L12").Value = 1252.5
$ws.Range("N12").Value = -1598.5

$ws.Range("H37").Value = 74232.57000000001
$ws.Range("J37").Value = 74232.57000000001
$ws.Range("L37").Value = 222697.71
$ws.Range("N37").Value = -222921.71

$ws.Range("H57").Value = 2500
$ws.Range("I57").Value = 2333.3333
$ws.Range("K57").Value = 6999.999899999999
$ws.Range("M57").Value = -6440.999899999999

$ws.Range("H132").Value = 46855.316
$ws.Range("I132").Value = 51190.85
$ws.Range("K132").Value = 460717.65
$ws.Range("M132").Value = -458187.65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4697.5
$ws.Range("I102").Value = 5763.3335
$ws.Range("K102").Value = 5763.3335
$ws.Range("M102").Value = -4141.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7678.84
$ws.Range("I46").Value = 1500.5
$ws.Range("J46").Value = 8216.087
$ws.Range("K46").Value = 1500.5
$ws.Range("L46").Value = 8216.087
$ws.Range("M46").Value = -1312.5
$ws.Range("N46").Value = -8592.087

$ws.Range("H68").Value = 2480.125
$ws.Range("I68").Value = 2274.6
$ws.Range("J68").Value = 2822.6667
$ws.Range("K68").Value = 2274.6
$ws.Range("L68").Value = 2822.6667
$ws.Range("M68").Value = -1525.6
$ws.Range("N68").Value = -4320.6667

$ws.Range("H71").Value = 2480.125
$ws.Range("I71").Value = 2274.6
$ws.Range("J71").Value = 2822.6667
$ws.Range("K71").Value = 11373
$ws.Range("L71").Value = 14113.3335
$ws.Range("M71").Value = -7629
$ws.Range("N71").Value = -21601.3335

$ws.Range("H122").Value = 2394.5
$ws.Range("J122").Value = 3699
$ws.Range("L122").Value = 11097
$ws.Range("N122").Value = -15997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 3656.8
$ws.Range("I3").Value = 8000
$ws.Range("J3").Value = 2571
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 2571
$ws.Range("M3").Value = -7886
$ws.Range("N3").Value = -2799

$ws.Range("H62").Value = 3999.5

$ws.Range("H65").Value = 3999.5

$ws.Range("H113").Value = 4186
$ws.Range("I113").Value = 3249.5
$ws.Range("J113").Value = 5122.5
$ws.Range("K113").Value = 9748.5
$ws.Range("L113").Value = 15367.5
$ws.Range("M113").Value = -7578.5
$ws.Range("N113").Value = -19707.5

$ws.Range("H126").Value = 4262.4287
$ws.Range("I126").Value = 4171.4346
$ws.Range("J126").Value = 4681
$ws.Range("K126").Value = 12514.3038
$ws.Range("L126").Value = 14043
$ws.Range("M126").Value = -10044.3038
$ws.Range("N126").Value = -18983

$ws.Range("H132").Value = 7053.4883
$ws.Range("I132").Value = 6368.6
$ws.Range("K132").Value = 19105.8
$ws.Range("M132").Value = -16575.8
